$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the width of columns D through M by 4 characters each.
$ws.Columns.Item(4).ColumnWidth = 34.17
$ws.Columns.Item(5).ColumnWidth = 36.17
$ws.Columns.Item(6).ColumnWidth = 24.17
$ws.Columns.Item(7).ColumnWidth = 26.17
$ws.Columns.Item(8).ColumnWidth = 21.17
$ws.Columns.Item(9).ColumnWidth = 23.17
$ws.Columns.Item(10).ColumnWidth = 29.17
$ws.Columns.Item(11).ColumnWidth = 31.17
$ws.Columns.Item(12).ColumnWidth = 24.17
$ws.Columns.Item(13).ColumnWidth = 26.17

# Rename the header labels from "link_testProject_*" to "link_project_*".
$ws.Range("D1").Value = "link_project_internalRoleLinkName"
$ws.Range("E1").Value = "link_project_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_project_project_id"
$ws.Range("G1").Value = "link_project_project_id_1"
$ws.Range("H1").Value = "link_project_team_id"
$ws.Range("I1").Value = "link_project_team_id_1"
$ws.Range("J1").Value = "link_project_test_project_id"
$ws.Range("K1").Value = "link_project_test_project_id_1"
$ws.Range("L1").Value = "link_project_trNthChild"
$ws.Range("M1").Value = "link_project_trNthChild_1"
